# Commit: "take care of one TODO"
#
# The generic "Custom AmountWithCurrency" component type used for the
# "Absolute Share" fields is replaced by the more specific
# "Custom EuTaxonomyAmountWithCurrencyComponent" component type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Framework Data Model")

$newComponent = "Custom EuTaxonomyAmountWithCurrencyComponent"

$targetCells = @("F12", "F14", "F16", "F19", "F31", "F33", "F35", "F38", "F50", "F52", "F54", "F57")
foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).Value = $newComponent
}

# Restore the active selection to where the author last left the cursor.
$null = $ws.Range("E15").Select()
